# Natmi following Dr Hou advice
# Update LR-pair statistics (Tgfb3-Tgfbr1) for rows 2-10 per updated cluster counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.825549
$ws.Range("H2").Value = 5.476647
$ws.Range("I2").Value = 0.04696949406168958
$ws.Range("J2").Value = 0.04696949406168958
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 63.76294533333333
$ws.Range("N2").Value = 191.288836
$ws.Range("O2").Value = 0.6446527016991613
$ws.Range("P2").Value = 0.6446527016991614
$ws.Range("Q2").Value = 116.4023810903213
$ws.Range("R2").Value = 1047.621429812892
$ws.Range("S2").Value = 0.03027901124431091
$ws.Range("T2").Value = 0.03027901124431091

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.825549
$ws.Range("H3").Value = 5.476647
$ws.Range("I3").Value = 0.04696949406168958
$ws.Range("J3").Value = 0.04696949406168958
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 22.32219333333333
$ws.Range("N3").Value = 66.96658
$ws.Range("O3").Value = 0.2256806388876402
$ws.Range("P3").Value = 0.2256806388876402
$ws.Range("Q3").Value = 40.75025771747332
$ws.Range("R3").Value = 366.7523194572599
$ws.Range("S3").Value = 0.01060010542807133
$ws.Range("T3").Value = 0.01060010542807133

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.825549
$ws.Range("H4").Value = 5.476647
$ws.Range("I4").Value = 0.04696949406168958
$ws.Range("J4").Value = 0.04696949406168958
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.825399
$ws.Range("N4").Value = 38.476197
$ws.Range("O4").Value = 0.1296666594131984
$ws.Range("P4").Value = 0.1296666594131984
$ws.Range("Q4").Value = 23.413394319051
$ws.Range("R4").Value = 210.720548871459
$ws.Range("S4").Value = 0.006090377389307349
$ws.Range("T4").Value = 0.006090377389307349

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 23.57737633333333
$ws.Range("H5").Value = 70.73212899999999
$ws.Range("I5").Value = 0.6066215903701957
$ws.Range("J5").Value = 0.6066215903701957
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 63.76294533333333
$ws.Range("N5").Value = 191.288836
$ws.Range("O5").Value = 0.6446527016991613
$ws.Range("P5").Value = 0.6446527016991614
$ws.Range("Q5").Value = 1503.36295824576
$ws.Range("R5").Value = 13530.26662421184
$ws.Range("S5").Value = 0.3910602471411886
$ws.Range("T5").Value = 0.3910602471411887

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.57737633333333
$ws.Range("H6").Value = 70.73212899999999
$ws.Range("I6").Value = 0.6066215903701957
$ws.Range("J6").Value = 0.6066215903701957
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 22.32219333333333
$ws.Range("N6").Value = 66.96658
$ws.Range("O6").Value = 0.2256806388876402
$ws.Range("P6").Value = 0.2256806388876402
$ws.Range("Q6").Value = 526.2987528054243
$ws.Range("R6").Value = 4736.688775248818
$ws.Range("S6").Value = 0.1369027480777821
$ws.Range("T6").Value = 0.1369027480777821

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 23.57737633333333
$ws.Range("H7").Value = 70.73212899999999
$ws.Range("I7").Value = 0.6066215903701957
$ws.Range("J7").Value = 0.6066215903701957
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.825399
$ws.Range("N7").Value = 38.476197
$ws.Range("O7").Value = 0.1296666594131984
$ws.Range("P7").Value = 0.1296666594131984
$ws.Range("Q7").Value = 302.3892588481569
$ws.Range("R7").Value = 2721.503329633412
$ws.Range("S7").Value = 0.07865859515122492
$ws.Range("T7").Value = 0.07865859515122492

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 13.46376966666667
$ws.Range("H8").Value = 40.391309
$ws.Range("I8").Value = 0.3464089155681148
$ws.Range("J8").Value = 0.3464089155681148
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 63.76294533333333
$ws.Range("N8").Value = 191.288836
$ws.Range("O8").Value = 0.6446527016991613
$ws.Range("P8").Value = 0.6446527016991614
$ws.Range("Q8").Value = 858.4896092362582
$ws.Range("R8").Value = 7726.406483126324
$ws.Range("S8").Value = 0.2233134433136618
$ws.Range("T8").Value = 0.2233134433136619

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 13.46376966666667
$ws.Range("H9").Value = 40.391309
$ws.Range("I9").Value = 0.3464089155681148
$ws.Range("J9").Value = 0.3464089155681148
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 22.32219333333333
$ws.Range("N9").Value = 66.96658
$ws.Range("O9").Value = 0.2256806388876402
$ws.Range("P9").Value = 0.2256806388876402
$ws.Range("Q9").Value = 300.5408694948022
$ws.Range("R9").Value = 2704.86782545322
$ws.Range("S9").Value = 0.07817778538178674
$ws.Range("T9").Value = 0.07817778538178674

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 13.46376966666667
$ws.Range("H10").Value = 40.391309
$ws.Range("I10").Value = 0.3464089155681148
$ws.Range("J10").Value = 0.3464089155681148
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.825399
$ws.Range("N10").Value = 38.476197
$ws.Range("O10").Value = 0.1296666594131984
$ws.Range("P10").Value = 0.1296666594131984
$ws.Range("Q10").Value = 172.678218019097
$ws.Range("R10").Value = 1554.103962171873
$ws.Range("S10").Value = 0.04491768687266614
$ws.Range("T10").Value = 0.04491768687266614

Write-Output "Updated rows 2-10 with new expression/specificity values"
